# Automatic update of the results tracker ("Actualizacion automatica del tracker")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while preventing Excel's automatic
# date-recognition from turning "YYYY-MM-DD" style strings into date serials.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# --- 1. Fill in the previously-pending result for row 72 -------------------
# G72 was blank ("") -> "Acierto"; H72 was blank ("") -> 0.73
$ws.Cells.Item(72, 7).Value = "Acierto"
$ws.Cells.Item(72, 8).Value = 0.73

# --- 2. Append the new match rows (76-79) -----------------------------------
$newRows = @(
    @{ Row = 76; A = 14762064; B = "2025-10-03"; C = "Francisco Cerundolo"; D = "Adrian Mannarino"; E = "Gana Francisco Cerundolo"; F = 1.8 },
    @{ Row = 77; A = 14762062; B = "2025-10-03"; C = "Jaume Munar"; D = "Flavio Cobolli"; E = "Gana Flavio Cobolli"; F = 1.8 },
    @{ Row = 78; A = 14786210; B = "2025-10-02"; C = "Joao Lucas Reis Da Silva"; D = "Facundo Bagnis"; E = "Gana Joao Lucas Reis Da Silva"; F = 2.1 },
    @{ Row = 79; A = 14763651; B = "2025-10-03"; C = "Viktorija Golubic"; D = "Alexandra Eala"; E = "Gana Alexandra Eala"; F = 2 }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    Set-TextCell $r 2 $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    # G/H (resultado/profit) stay blank/pending, just like the source rows
    # that are awaiting a result -- leave them unset.
}
